$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "49.772.86"
$ws.Range("E2").Value = "  -0.30%  "

$ws.Range("D3").Value = "2.651.50"
$ws.Range("E3").Value = "  +0.29%  "

$ws.Range("E4").Value = "  +0.00%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "113.07"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.49%  "

$ws.Range("E6").Value = "  +0.45%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.525"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.84%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.999"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -0.05%  "

$ws.Range("E9").Value = "  -0.80%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "39.88"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -2.55%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "20.07"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.03%  "

$ws.Range("E12").Value = "  -0.58%  "

$ws.Range("E13").Value = "  +2.18%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "7.59"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +2.87%  "

$ws.Range("D15").Value = "3.065.92"
$ws.Range("E15").Value = "  +0.23%  "

$ws.Range("D16").Value = "2.653.14"
$ws.Range("E16").Value = "  +0.19%  "

$ws.Range("E17").Value = "  -1.09%  "

$ws.Range("D18").Value = "49.737.63"
$ws.Range("E18").Value = "  -0.27%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "13.38"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +1.46%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "2.94"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.65%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.71"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.57%  "

$ws.Range("E22").Value = "  -0.40%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "269.24"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -2.23%  "

$ws.Range("E24").Value = "  -4.00%  "

$ws.Range("E25").Value = "  -0.69%  "

$ws.Range("E26").Value = "  -2.29%  "

$ws.Range("E27").Value = "  +0.01%  "

$ws.Range("E28").Value = "  +1.97%  "

$ws.Range("E29").Value = "  -0.66%  "

$ws.Range("B30").Value = "InjectiveProtocol"
$ws.Range("C30").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "35.03"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -3.75%  "

$ws.Range("B31").Value = "Kaspa"
$ws.Range("C31").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.138"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -2.12%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.51"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.95%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.0820"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +0.33%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "19.23"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -1.19%  "

$ws.Range("E36").Value = "  -0.16%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "4.96"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -1.33%  "

$ws.Range("E38").Value = "  -1.13%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "3.14"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +1.10%  "

$ws.Range("B40").Value = "EnergySwap"
$ws.Range("C40").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "23.79"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +8.18%  "

$ws.Range("B41").Value = "Monero"
$ws.Range("C41").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "128.94"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +4.01%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.0344"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +8.69%  "

$ws.Range("E43").Value = "  +2.27%  "

$ws.Range("E44").Value = "  -0.41%  "

$ws.Range("E45").Value = "  +0.07%  "

$ws.Range("D46").Value = "2.068.72"
$ws.Range("E46").Value = "  -0.88%  "

$ws.Range("E47").Value = "  +7.60%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.20"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -2.91%  "

$ws.Range("E49").Value = "  -2.21%  "

$ws.Range("E50").Value = "  -1.29%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "59.32"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.70%  "
